$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; this shifts the existing rows 11-42 down to 12-43
$ws.Rows.Item(11).Insert()

# Fill the newly inserted row 11 with the new record's data.
# Columns A, B, C, E, F, G, I, N, Q, R are identical to the surrounding rows.
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 44533
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100112031
$ws.Cells.Item(11, 7).Value = "Poroto verde"
$ws.Cells.Item(11, 8).Value = "Magnum"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 29000
$ws.Cells.Item(11, 12).Value = 30000
$ws.Cells.Item(11, 13).Value = 29500
$ws.Cells.Item(11, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 1180
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"
